# The workbook was previously authored/saved by a German-locale Excel,
# which gave the single worksheet its localized default name "Tabelle1".
# The project now loads/writes this correction data via pandas (see the
# commit message), whose default sheet name convention is the English
# "Sheet1" - so rename the sheet accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"
